$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.441.23"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "'2.504.40"
$ws.Range("E3").Value = "  -4.94%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'581.47"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'171.54"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "'2.504.04"
$ws.Range("E9").Value = "  -4.93%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'26.61"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "'2.957.13"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").Value = "'66.329.61"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "'2.493.64"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("E19").Value = "  -6.25%  "
$ws.Range("D20").Value = "'7.66"
$ws.Range("E20").Value = "  -4.59%  "
$ws.Range("D21").Value = "'347.06"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "'4.61"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D26").Value = "'69.52"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'2.626.89"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("D30").Value = "'0.0₃0975"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "'527.33"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").Value = "'156.29"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'18.59"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'5.08"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").Value = "'39.48"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "'0.556"
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("D49").Value = "'3.67"
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "'0.0₆0272"
$ws.Range("E50").Value = "  -9.36%  "
$ws.Range("E51").Value = "  +0.79%  "
